# Apply weekly Fruta/hortaliza update:
# Insert two new price-report rows at row 95 (pushing the existing rows 95-209
# down to 97-211) and populate them with the new "Artic Pride" / "Big John"
# entries dated 2021-12-09 (Excel serial 44539).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 95.
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()

# --- New row 95: Artic Pride / Primera ---
$ws.Range('A95').Value() = 11
$ws.Range('B95').Value() = "Vega Monumental Concepción"
$ws.Range('C95').Value() = "Bíobío"
$ws.Range('D95').Value() = 44539
$ws.Range('E95').Value() = 8
$ws.Range('F95').Value() = "Fruta"
$ws.Range('G95').Value() = 100103
$ws.Range('H95').Value() = "Frutos de hueso (carozo)"
$ws.Range('I95').Value() = 100103006
$ws.Range('J95').Value() = "Nectarín"
$ws.Range('K95').Value() = "Artic Pride"
$ws.Range('L95').Value() = "Primera"
$ws.Range('M95').Value() = 220
$ws.Range('N95').Value() = 13000
$ws.Range('O95').Value() = 14000
$ws.Range('P95').Value() = 13545
$ws.Range('Q95').Value() = "`$/caja 15 kilos empedrada"
$ws.Range('R95').Value() = "Región de O'Higgins"
$ws.Range('S95').Value() = 903
$ws.Range('T95').Value() = 15

# --- New row 96: Big John / Primera ---
$ws.Range('A96').Value() = 11
$ws.Range('B96').Value() = "Vega Monumental Concepción"
$ws.Range('C96').Value() = "Bíobío"
$ws.Range('D96').Value() = 44539
$ws.Range('E96').Value() = 8
$ws.Range('F96').Value() = "Fruta"
$ws.Range('G96').Value() = 100103
$ws.Range('H96').Value() = "Frutos de hueso (carozo)"
$ws.Range('I96').Value() = 100103006
$ws.Range('J96').Value() = "Nectarín"
$ws.Range('K96').Value() = "Big John"
$ws.Range('L96').Value() = "Primera"
$ws.Range('M96').Value() = 220
$ws.Range('N96').Value() = 12000
$ws.Range('O96').Value() = 13000
$ws.Range('P96').Value() = 12455
$ws.Range('Q96').Value() = "`$/caja 15 kilos empedrada"
$ws.Range('R96').Value() = "Región de O'Higgins"
$ws.Range('S96').Value() = 830
$ws.Range('T96').Value() = 15

# Make sure the date cells use the existing date number format, same
# as the rest of column D (the Insert() above already carried the style
# down, but set it explicitly here for robustness).
$dateFormat = $ws.Range('D97').NumberFormat
$ws.Range('D95').NumberFormat = $dateFormat
$ws.Range('D96').NumberFormat = $dateFormat
